$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("G2").Value = 201.4397426666667
$ws.Range("H2").Value = 604.3192280000001
$ws.Range("I2").Value = 0.4833500233086392
$ws.Range("J2").Value = 0.4833500233086393
$ws.Range("M2").Value = 0.07215833333333334
$ws.Range("N2").Value = 0.216475
$ws.Range("Q2").Value = 14.53555609792222
$ws.Range("R2").Value = 130.8200048813
$ws.Range("S2").Value = 0.4833500233086392
$ws.Range("T2").Value = 0.4833500233086393

# Row 3
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("I3").Value = 0.1569674599353791
$ws.Range("J3").Value = 0.1569674599353792
$ws.Range("M3").Value = 0.07215833333333334
$ws.Range("N3").Value = 0.216475
$ws.Range("Q3").Value = 4.720408005405556
$ws.Range("R3").Value = 42.48367204865
$ws.Range("S3").Value = 0.1569674599353791
$ws.Range("T3").Value = 0.1569674599353792

# Row 4
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 60.43484133333334
$ws.Range("H4").Value = 181.304524
$ws.Range("I4").Value = 0.1450120099461104
$ws.Range("J4").Value = 0.1450120099461104
$ws.Range("M4").Value = 0.07215833333333334
$ws.Range("N4").Value = 0.216475
$ws.Range("Q4").Value = 4.360877425877779
$ws.Range("R4").Value = 39.2478968329
$ws.Range("S4").Value = 0.1450120099461104
$ws.Range("T4").Value = 0.1450120099461104

# Row 5
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 89.46554166666668
$ws.Range("H5").Value = 268.396625
$ws.Range("I5").Value = 0.2146705068098712
$ws.Range("J5").Value = 0.2146705068098712
$ws.Range("M5").Value = 0.07215833333333334
$ws.Range("N5").Value = 0.216475
$ws.Range("Q5").Value = 6.455684377430557
$ws.Range("R5").Value = 58.101159396875
$ws.Range("S5").Value = 0.2146705068098712
$ws.Range("T5").Value = 0.2146705068098712
